$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.73%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.77%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.131"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.32%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08171"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.88%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.983"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'2.09%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'2.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9376"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.00%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1299"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-7.51%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.98%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-0.62%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'0.16%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09731"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.79%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001408"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.25%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.006057"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.65%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.633"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-7.87%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.370"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'3.36%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.280"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.91%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3491"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.84%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-1.58%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.968"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.60%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2580"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'6.53%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04357"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.73%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.99%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004763"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'9.50%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'199.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-7.59%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02213"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'8.97%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05201"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.90%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.64%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01037"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.11%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1399"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.68%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-1.42%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'6.83%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006945"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'9.36%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002883"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.67%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'30.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
